# Apply the updated "Price" (D) and "Volume(1h)" (E) figures for the
# cryptocurrency rows, as refreshed by the scheduled GitHub Actions run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.495.01"
$ws.Range("E2").Value = "  +2.22%  "
$ws.Range("D3").Value = "'1.841.16"
$ws.Range("E3").Value = "  +1.59%  "
$ws.Range("E4").Value = "  +1.17%  "
$ws.Range("D5").Value = "'314.86"
$ws.Range("E5").Value = "  +1.83%  "
$ws.Range("D7").Value = "'0.4753"
$ws.Range("E7").Value = "  +1.96%  "
$ws.Range("D8").Value = "'0.3698"
$ws.Range("E8").Value = "  +1.25%  "
$ws.Range("E9").Value = "  +1.90%  "
$ws.Range("D10").Value = "'0.8864"
$ws.Range("D11").Value = "'20.48"
$ws.Range("D12").Value = "'1.880.07"
$ws.Range("E12").Value = "  +4.82%  "
$ws.Range("D13").Value = "'0.07361"
$ws.Range("E13").Value = "  +4.00%  "
$ws.Range("D14").Value = "'5.457"
$ws.Range("E14").Value = "  +1.46%  "
$ws.Range("D15").Value = "'93.30"
$ws.Range("E15").Value = "  +1.88%  "
$ws.Range("E16").Value = "  +1.39%  "
$ws.Range("E17").Value = "  +0.99%  "
$ws.Range("D18").Value = "'0.000008829"
$ws.Range("E18").Value = "  +1.59%  "
$ws.Range("D19").Value = "'1.013"
$ws.Range("E19").Value = "  +1.10%  "
$ws.Range("D20").Value = "'14.83"
$ws.Range("E20").Value = "  +1.50%  "
$ws.Range("D21").Value = "'27.501.84"
$ws.Range("E21").Value = "  +2.19%  "
$ws.Range("D22").Value = "'5.331"
$ws.Range("E22").Value = "  +0.78%  "
$ws.Range("E23").Value = "  +0.98%  "
$ws.Range("D24").Value = "'2.101.45"
$ws.Range("E24").Value = "  +3.43%  "
$ws.Range("D25").Value = "'1.901"
$ws.Range("E25").Value = "  +0.35%  "
$ws.Range("D26").Value = "'152.29"
$ws.Range("E26").Value = "  +1.21%  "
$ws.Range("D27").Value = "'18.64"
$ws.Range("E27").Value = "  +2.16%  "
$ws.Range("D28").Value = "'2.157"
$ws.Range("E28").Value = "  +0.22%  "
$ws.Range("D29").Value = "'5.251"
$ws.Range("E29").Value = "  -0.28%  "
$ws.Range("D30").Value = "'118.16"
$ws.Range("E30").Value = "  +2.47%  "
$ws.Range("D31").Value = "'0.09009"
$ws.Range("E31").Value = "  +0.87%  "
$ws.Range("D32").Value = "'0.7583"
$ws.Range("E32").Value = "  +0.73%  "
$ws.Range("E33").Value = "  +2.42%  "
$ws.Range("D34").Value = "'4.563"
$ws.Range("E34").Value = "  +1.82%  "
$ws.Range("D35").Value = "'2.968"
$ws.Range("E35").Value = "  +1.92%  "
$ws.Range("E36").Value = "  +1.23%  "
$ws.Range("E37").Value = "  +2.30%  "
$ws.Range("D38").Value = "'0.05346"
$ws.Range("E38").Value = "  +1.29%  "
$ws.Range("D39").Value = "'0.01961"
$ws.Range("E39").Value = "  +0.80%  "
$ws.Range("D40").Value = "'3.005"
$ws.Range("E40").Value = "  +0.61%  "
$ws.Range("D41").Value = "'7.341"
$ws.Range("E41").Value = "  +2.09%  "
$ws.Range("D42").Value = "'2.403"
$ws.Range("E42").Value = "  +5.33%  "
$ws.Range("D43").Value = "'0.5344"
$ws.Range("E43").Value = "  +1.05%  "
$ws.Range("E44").Value = "  +0.85%  "
$ws.Range("D45").Value = "'8.540"
$ws.Range("E45").Value = "  +1.88%  "
$ws.Range("D46").Value = "'0.4927"
$ws.Range("E46").Value = "  +1.44%  "
$ws.Range("D47").Value = "'10.55"
$ws.Range("E47").Value = "  +1.43%  "
$ws.Range("D48").Value = "'1.014"
$ws.Range("E48").Value = "  +1.21%  "
$ws.Range("D49").Value = "'104.86"
$ws.Range("E49").Value = "  +1.81%  "
$ws.Range("E50").Value = "  +1.60%  "
$ws.Range("D51").Value = "'0.06323"
$ws.Range("E51").Value = "  +0.51%  "
